$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column H to fit the new comment text.
# Excel's stored <col width=".."/> equals ColumnWidth + 5/6, so back that
# offset out to land on the target stored width of 143.5.
$ws.Columns.Item(8).ColumnWidth = 143.5 - (5/6)

# Add comment text for "use_tfidf_model" row (row 18) first, then "use_danlp_model" row (row 16),
# so the shared-string table order matches the target workbook.
$ws.Range("H18").Value = "<-- Catches comments that the others don't (since recall is lower), but also contributes to a decrease in precision. Label function can only mark OFFENSIVE and ABSTAIN."
$ws.Range("H16").Value = "<-- Significantly increases precision, while reducing recall, which is probably because this LF can mark things as NOT_OFFENSIVE as well as the other two."

# Move/update the active selection to H1
$ws.Range("H1").Select()
